$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column J (2021 data) added next to existing I (2020) column.
# For each data row, copy the formatting from the corresponding H-column
# cell (the engine reused those style indices when the workbook was
# authored) and then set the new 2021 value.

$ws.Range("J4").Value = 2021
$ws.Range("H4").Copy()
$ws.Range("J4").PasteSpecial(-4122)

$ws.Range("J5").Value = 1006091.2
$ws.Range("H5").Copy()
$ws.Range("J5").PasteSpecial(-4122)

$ws.Range("J6").Value = 2092.6999999999998
$ws.Range("H6").Copy()
$ws.Range("J6").PasteSpecial(-4122)

$ws.Range("J7").Value = 211904.6
$ws.Range("H7").Copy()
$ws.Range("J7").PasteSpecial(-4122)

$ws.Range("J8").Value = 228945.8
$ws.Range("H8").Copy()
$ws.Range("J8").PasteSpecial(-4122)

$ws.Range("J9").Value = 6780.6
$ws.Range("H9").Copy()
$ws.Range("J9").PasteSpecial(-4122)

$ws.Range("J10").Value = 92.5
$ws.Range("H10").Copy()
$ws.Range("J10").PasteSpecial(-4122)

$ws.Range("J11").Value = 9456.7999999999993
$ws.Range("H11").Copy()
$ws.Range("J11").PasteSpecial(-4122)

$ws.Range("J12").Value = 92470.9
$ws.Range("H12").Copy()
$ws.Range("J12").PasteSpecial(-4122)

$ws.Range("J13").Value = 656.4
$ws.Range("H13").Copy()
$ws.Range("J13").PasteSpecial(-4122)

$ws.Range("J14").Value = 3692
$ws.Range("H14").Copy()
$ws.Range("J14").PasteSpecial(-4122)

$ws.Range("J15").Value = 59559.1
$ws.Range("H15").Copy()
$ws.Range("J15").PasteSpecial(-4122)

$ws.Range("J16").Value = 53592.2
$ws.Range("H16").Copy()
$ws.Range("J16").PasteSpecial(-4122)

$ws.Range("J17").Value = 11799.2
$ws.Range("H17").Copy()
$ws.Range("J17").PasteSpecial(-4122)

$ws.Range("J18").Value = 316755
$ws.Range("H18").Copy()
$ws.Range("J18").PasteSpecial(-4122)

$ws.Range("J19").Value = 901
$ws.Range("H19").Copy()
$ws.Range("J19").PasteSpecial(-4122)

$ws.Range("J20").Value = 76.5
$ws.Range("H20").Copy()
$ws.Range("J20").PasteSpecial(-4122)

$ws.Range("J21").Value = 1672.3
$ws.Range("H21").Copy()
$ws.Range("J21").PasteSpecial(-4122)

# Row 22: H22's format uses the General number format, but J22 needs a
# thousands-separated decimal format, so a new cell style (based on H22's
# style) is created with that number format applied.
$ws.Range("J22").Value = 5539.9
$ws.Range("H22").Copy()
$ws.Range("J22").PasteSpecial(-4122)
$ws.Range("J22").NumberFormat = "#,##0.0"

$ws.Range("J23").Value = 103.7
$ws.Range("H23").Copy()
$ws.Range("J23").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Header row is now shorter since it wraps less with the extra column.
$ws.Rows.Item(1).RowHeight = 33.75

# Update selected cell shown when the sheet is reopened.
$ws.Range("K3").Select()
